$d = $word.ActiveDocument

function Find-ParaByPrefix($doc, [string]$prefix) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r").StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

function Find-ParaByExact($doc, [string]$value) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r") -eq $value) {
            return $p
        }
    }
    return $null
}

# --- Class Project ------------------------------------------------------
# Drop the "String []groupdID;...danh sach..." line entirely - the field
# is removed from the Project class.
$pGroupdID = Find-ParaByPrefix $d "String []groupdID;"
if ($pGroupdID -ne $null) { $pGroupdID.Range.Delete() }

# "Long ID;" + "String projectID;" (two separate fields right after the
# "// Admin ... projectmanager" comment inside Class Project) collapse
# into the single field "Long projectID;": drop the second paragraph and
# retarget the first one's text, scoping Find to that single paragraph so
# none of the many other "Long ID;" lines elsewhere get touched.
$pProjectID = Find-ParaByExact $d "String projectID;"
if ($pProjectID -ne $null) { $pProjectID.Range.Delete() }

$pComment = Find-ParaByPrefix $d "// Admin"
if ($pComment -ne $null) {
    $pLongID = $pComment.Next()
    $null = $pLongID.Range.Find.Execute("Long ID;", $true, $false, $false, $false, $false,
                             $true, 1, $false, "Long projectID;", 2)
}

# --- Class Group ---------------------------------------------------------
# "Boolean Work; " -> "Boolean isWorking; " (field renamed).
$pWork = Find-ParaByPrefix $d "Boolean Work;"
if ($pWork -ne $null) {
    $null = $pWork.Range.Find.Execute("Boolean Work;", $true, $false, $false, $false, $false,
                             $true, 1, $false, "Boolean isWorking;", 2)
}

# Drop the "String  idusers[];" line entirely - the field is removed from
# the Group class.
$pIdusers = Find-ParaByPrefix $d "String  idusers"
if ($pIdusers -ne $null) { $pIdusers.Range.Delete() }
